# Quarterly indexing esoteric bug-fix operation
#
# The dates in column A (rows 2-63) were being stamped with the 1st of the
# quarter-start month, one quarter too early. The fix re-derives each date
# as the 15th of the month *after* the originally stored month, preserving
# the year roll-over at December -> January.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 63 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($serial -eq $null) { continue }

    $oldDate = [DateTime]::FromOADate([double]$serial)

    $newMonth = $oldDate.Month + 1
    $newYear = $oldDate.Year
    if ($newMonth -gt 12) {
        $newMonth -= 12
        $newYear += 1
    }

    $newDate = Get-Date -Year $newYear -Month $newMonth -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value2 = [int][Math]::Floor($newDate.ToOADate())
}
